$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36 (existing rows 36-106 shift down to 37-107,
# carrying their formatting/styles with them).
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new record.
$ws.Range("A36").Value = 10
$ws.Range("B36").Value = "Vega Modelo de Temuco"
$ws.Range("C36").Value = "La Araucanía"
$ws.Range("D36").Value = 45259
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = 100112026
$ws.Range("G36").Value = "Haba"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 25
$ws.Range("K36").Value = 10000
$ws.Range("L36").Value = 10000
$ws.Range("M36").Value = 10000
$ws.Range("N36").Value = "$/saco 25 kilos"
$ws.Range("O36").Value = "Región del Maule"
$ws.Range("P36").Value = 400
$ws.Range("Q36").Value = 25
$ws.Range("R36").Value = "Hortaliza"
